$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.171.68'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.43%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.537.23'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.44%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '183.57'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.02%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.599'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.84%  '

$ws.Range("E9").Value = '  +4.12%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.15'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.82%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.441'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.68%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.142.36'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.36%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '32.78'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +13.21%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.134'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.12%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '68.129.51'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.40%  '

$ws.Range("E16").Value = '  +1.04%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.530.47'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.36%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.44'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.51%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.89'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.24%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '399.94'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.84%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.13'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.72%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.70'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.39%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.548'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.24%  '

$ws.Range("E24").Value = '  +0.21%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000127'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.51%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.69'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.01%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.73'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.14%  '

$ws.Range("E28").Value = '  -0.69%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.997'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.37'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.56%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.49'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.76%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.09'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.77%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '24.17'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.55%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.53'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.39%  '

$ws.Range("E35").Value = '  +0.08%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.69'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.27%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '164.40'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.30%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.98'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.07%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.883'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.43%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.17'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.56%  '

$ws.Range("E41").Value = '  +7.54%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.78'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.75%  '

$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '27.21'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.61%  '

$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '27.91'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.74%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.882.03'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.84%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0746'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.21%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '42.43'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.10%  '

$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '352.30'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.75%  '

$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0309'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.99%  '

$ws.Range("E50").Value = '  +0.52%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '34.23'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.14%  '
